$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------------
# Paragraph (containing the hyperlinked equation refs "1, 2, 3, 4") used to
# read: "... infectées, " recovered" (trouver traduction) Le modèle ...".
# The undefined "recovered"/"(trouver traduction)" placeholder is replaced
# with the real French term "immunisées." so the sentence reads:
# "... infectées, immunisées. Le modèle ...".
$find1 = "faisant référence respectivement au nombre de personnes sensisbles, exposées, infectées, “recovered” (trouver traduction) Le modèle de dynamique épidémiologique pour le moustique est celui-ci :"
$replace1 = "faisant référence respectivement au nombre de personnes sensisbles, exposées, infectées, immunisées. Le modèle de dynamique épidémiologique pour le moustique est celui-ci :"

$found1 = $d.Content.Find.Execute($find1, $true, $false, $false, $false, $false, $true, 1, $false, $replace1, 2)
if (-not $found1) {
    throw "Change 1: target text not found"
}

# --- Change 2 ---------------------------------------------------------
# "... ont modifié ce troisième modèle pour forcer certaines variables ..."
# becomes "... pour influencer certaines variables ...".
$find2 = "Dans un troisième temps, les auteurs de cet article ont modifié ce troisième modèle pour forcer certaines variables controllant le modèle 2. Les variables forcés par la température sont :"
$replace2 = "Dans un troisième temps, les auteurs de cet article ont modifié ce troisième modèle pour influencer certaines variables controllant le modèle 2. Les variables forcés par la température sont :"

$found2 = $d.Content.Find.Execute($find2, $true, $false, $false, $false, $false, $true, 1, $false, $replace2, 2)
if (-not $found2) {
    throw "Change 2: target text not found"
}
